$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Packet Calculations" (sheet1)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Packet Calculations")

# Insert three new rows:
#   - row 11 -> new "Temp" sensor row   (old row11 "Mag" shifts to 12, etc.)
#   - row 13 (after first insert) -> new "Time" sensor row
#   - row 19 (after first two inserts) -> new "Size" header row
$ws1.Rows.Item(11).Insert()
$ws1.Rows.Item(13).Insert()
$ws1.Rows.Item(19).Insert()

# New "Temp" row (row 11)
$ws1.Range("A11").Value = "Temp"
$ws1.Range("B11").Value = 1
$ws1.Range("C11").Value = 16
$ws1.Range("D11").Value = 2
$ws1.Range("E11").Formula = "=B11*D11"

# "Mag" row, shifted down to row 12 - make sure formulas/values are intact
$ws1.Range("A12").Value = "Mag"
$ws1.Range("B12").Value = 3
$ws1.Range("C12").Value = 16
$ws1.Range("D12").Formula = "=C12/8"
$ws1.Range("E12").Formula = "=B12*D12"

# New "Time" row (row 13)
$ws1.Range("A13").Value = "Time "
$ws1.Range("B13").Value = 1
$ws1.Range("C13").Value = 16
$ws1.Range("D13").Formula = "=C13/8"
$ws1.Range("E13").Formula = "=B13*D13"

# "Total Reading Size" row, shifted to row 14
$ws1.Range("E14").Formula = "=SUM(E8:E13)"

# "Header" row shifted to row 17 (label only, no change needed besides shift)
# "Data" / Size (bits) / Size (bytes) header row shifted to row 18 (no change needed)

# New "Size" row (row 19)
$ws1.Range("A19").Value = "Size"
$ws1.Range("B19").Value = 16
$ws1.Range("C19").Value = 2

# "Counter" row shifted to row 20
$ws1.Range("A20").Value = "Counter"
$ws1.Range("B20").Value = 16
$ws1.Range("C20").Formula = "=B20/8"

# "Total Header size" row shifted to row 21
$ws1.Range("C21").Formula = "=SUM(C19:C20)"

# Update the references at the top of the sheet (I4:J6) so they point to the
# new "Total Header size" (C21) / "Total Reading Size" (E14) cells.
$ws1.Range("J4").Formula = "=B1-C21"
$ws1.Range("J5").Formula = "=_xlfn.FLOOR.MATH(J4/E14)"
$ws1.Range("J6").Formula = "=C21+J5*E14"
$ws1.Range("B2").Formula = "=J6"

# ---------------------------------------------------------------------------
# Sheet "Timing" (sheet2)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Timing")

# New columns for the cost calculations
$ws2.Range("I4").Value = "Cost per packet"
$ws2.Range("J4").Value = "Cost for all data"
$ws2.Range("K4").Value = "Cost with resending "

$ws2.Range("I5").Formula = "=0.04 + 0.0015 * ('Packet Calculations'!J6-'Packet Calculations'!D23-30)"
$ws2.Range("J5").Formula = "=I5*F5"
$ws2.Range("K5").Formula = "=J5*3"

$ws2.Columns.Item(9).ColumnWidth = 12.3
$ws2.Columns.Item(10).ColumnWidth = 16.65
$ws2.Columns.Item(11).ColumnWidth = 12.0

# Estimated flight time changed from 200 to 350
$ws2.Range("F2").Value = 350

# ---------------------------------------------------------------------------
# Selections (match committed cursor positions)
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("B21").Select()

$ws2.Activate()
$ws2.Range("H8").Select()
